$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the missing "Avg Accuracy" value for the last row (I11)
$ws.Range("I11").Value = 0.618319801057855

# Update the active selection to B12 (next row below the table), matching the saved selection state
$ws.Range("B12").Select()
